# Auto-generated: update cached market-data values per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 100.625
$ws.Cells.Item(9, 9).Value = 97.85714
$ws.Cells.Item(9, 11).Value = 97.85714
$ws.Cells.Item(9, 13).Value = 71.14286
# Row 18
$ws.Cells.Item(18, 8).Value = 12649.292
$ws.Cells.Item(18, 9).Value = 8963.637000000001
$ws.Cells.Item(18, 10).Value = 15767.923
$ws.Cells.Item(18, 11).Value = 8963.637000000001
$ws.Cells.Item(18, 12).Value = 15767.923
$ws.Cells.Item(18, 13).Value = -8679.637000000001
$ws.Cells.Item(18, 14).Value = -16335.923
# Row 32
$ws.Cells.Item(32, 8).Value = 1533.4445
$ws.Cells.Item(32, 9).Value = 267
$ws.Cells.Item(32, 11).Value = 267
$ws.Cells.Item(32, 13).Value = 59
# Row 70
$ws.Cells.Item(70, 8).Value = 7605.515
$ws.Cells.Item(70, 9).Value = 5230.375
$ws.Cells.Item(70, 10).Value = 8365.559999999999
$ws.Cells.Item(70, 11).Value = 15691.125
$ws.Cells.Item(70, 12).Value = 25096.68
$ws.Cells.Item(70, 13).Value = -15421.125
$ws.Cells.Item(70, 14).Value = -25636.68
# Row 73
$ws.Cells.Item(73, 8).Value = 7605.515
$ws.Cells.Item(73, 9).Value = 5230.375
$ws.Cells.Item(73, 10).Value = 8365.559999999999
$ws.Cells.Item(73, 11).Value = 15691.125
$ws.Cells.Item(73, 12).Value = 25096.68
$ws.Cells.Item(73, 13).Value = -14755.125
$ws.Cells.Item(73, 14).Value = -26968.68
# Row 131
$ws.Cells.Item(131, 8).Value = 3114.6843
$ws.Cells.Item(131, 10).Value = 4020
$ws.Cells.Item(131, 12).Value = 12060
$ws.Cells.Item(131, 14).Value = -22140
# Row 132
$ws.Cells.Item(132, 8).Value = 853.3946999999999
$ws.Cells.Item(132, 9).Value = 839.7222
$ws.Cells.Item(132, 11).Value = 2519.1666
$ws.Cells.Item(132, 13).Value = 10.83339999999998
# Row 133
$ws.Cells.Item(133, 8).Value = 70000
$ws.Cells.Item(133, 10).Value = 70000
$ws.Cells.Item(133, 12).Value = 70000
$ws.Cells.Item(133, 14).Value = -80120
# Row 137
$ws.Cells.Item(137, 8).Value = 1708.7307
$ws.Cells.Item(137, 9).Value = 1429.9333
$ws.Cells.Item(137, 11).Value = 4289.7999
$ws.Cells.Item(137, 13).Value = -1739.7999
# Row 138
$ws.Cells.Item(138, 8).Value = 3224.158
$ws.Cells.Item(138, 9).Value = 5525.636
$ws.Cells.Item(138, 10).Value = 2286.5186
$ws.Cells.Item(138, 11).Value = 16576.908
$ws.Cells.Item(138, 12).Value = 6859.5558
$ws.Cells.Item(138, 13).Value = -11436.908
$ws.Cells.Item(138, 14).Value = -17139.5558

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 4234.3774
$ws.Cells.Item(32, 9).Value = 3660.4878
$ws.Cells.Item(32, 10).Value = 6195.1665
$ws.Cells.Item(32, 11).Value = 3660.4878
$ws.Cells.Item(32, 12).Value = 6195.1665
$ws.Cells.Item(32, 13).Value = -3373.4878
$ws.Cells.Item(32, 14).Value = -6769.1665
# Row 122
$ws.Cells.Item(122, 8).Value = 1150
$ws.Cells.Item(122, 9).Value = 1000
$ws.Cells.Item(122, 10).Value = 1300
$ws.Cells.Item(122, 11).Value = 3000
$ws.Cells.Item(122, 12).Value = 3900
$ws.Cells.Item(122, 13).Value = -550
$ws.Cells.Item(122, 14).Value = -8800
# Row 132
$ws.Cells.Item(132, 8).Value = 1443.8572
$ws.Cells.Item(132, 9).Value = 1139.027
$ws.Cells.Item(132, 11).Value = 3417.081
$ws.Cells.Item(132, 13).Value = -887.0810000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 2267.4546
$ws.Cells.Item(20, 9).Value = 2316
$ws.Cells.Item(20, 11).Value = 2316
$ws.Cells.Item(20, 13).Value = -2069
# Row 86
$ws.Cells.Item(86, 8).Value = 2376.4443
$ws.Cells.Item(86, 9).Value = 2397.6
$ws.Cells.Item(86, 10).Value = 2350
$ws.Cells.Item(86, 11).Value = 2397.6
$ws.Cells.Item(86, 12).Value = 2350
$ws.Cells.Item(86, 13).Value = -1274.6
$ws.Cells.Item(86, 14).Value = -4596
# Row 89
$ws.Cells.Item(89, 8).Value = 2376.4443
$ws.Cells.Item(89, 9).Value = 2397.6
$ws.Cells.Item(89, 10).Value = 2350
$ws.Cells.Item(89, 11).Value = 11988
$ws.Cells.Item(89, 12).Value = 11750
$ws.Cells.Item(89, 13).Value = -6372
$ws.Cells.Item(89, 14).Value = -22982
# Row 107
$ws.Cells.Item(107, 8).Value = 3311.2222
$ws.Cells.Item(107, 9).Value = 3206.6875
$ws.Cells.Item(107, 11).Value = 3206.6875
$ws.Cells.Item(107, 13).Value = -1286.6875
# Row 129
$ws.Cells.Item(129, 8).Value = 44999.25
# Row 131
$ws.Cells.Item(131, 8).Value = 46629.5
$ws.Cells.Item(131, 10).Value = 46629.5
$ws.Cells.Item(131, 12).Value = 46629.5
$ws.Cells.Item(131, 14).Value = -56709.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 4167500.2
$ws.Cells.Item(22, 9).Value = 641.25
$ws.Cells.Item(22, 10).Value = 8929625
$ws.Cells.Item(22, 11).Value = 641.25
$ws.Cells.Item(22, 12).Value = 8929625
$ws.Cells.Item(22, 13).Value = -291.25
$ws.Cells.Item(22, 14).Value = -8930325
# Row 31
$ws.Cells.Item(31, 8).Value = 1370.6316
$ws.Cells.Item(31, 9).Value = 679.5454999999999
$ws.Cells.Item(31, 10).Value = 1805.0286
$ws.Cells.Item(31, 11).Value = 679.5454999999999
$ws.Cells.Item(31, 12).Value = 1805.0286
$ws.Cells.Item(31, 13).Value = -384.5454999999999
$ws.Cells.Item(31, 14).Value = -2395.0286
# Row 34
$ws.Cells.Item(34, 8).Value = 1370.6316
$ws.Cells.Item(34, 9).Value = 679.5454999999999
$ws.Cells.Item(34, 10).Value = 1805.0286
$ws.Cells.Item(34, 11).Value = 679.5454999999999
$ws.Cells.Item(34, 12).Value = 1805.0286
$ws.Cells.Item(34, 13).Value = -477.5454999999999
$ws.Cells.Item(34, 14).Value = -2209.0286
# Row 58
$ws.Cells.Item(58, 8).Value = 2899999.8
$ws.Cells.Item(58, 9).Value = 4832341
$ws.Cells.Item(58, 11).Value = 4832341
$ws.Cells.Item(58, 13).Value = -4832138
# Row 132
$ws.Cells.Item(132, 8).Value = 2196.9756
$ws.Cells.Item(132, 9).Value = 1598.4688
$ws.Cells.Item(132, 10).Value = 4325
$ws.Cells.Item(132, 11).Value = 4795.4064
$ws.Cells.Item(132, 12).Value = 12975
$ws.Cells.Item(132, 13).Value = -2265.4064
$ws.Cells.Item(132, 14).Value = -18035
# Row 134
$ws.Cells.Item(134, 8).Value = 1500.62
$ws.Cells.Item(134, 9).Value = 831.4872
$ws.Cells.Item(134, 11).Value = 2494.4616
$ws.Cells.Item(134, 13).Value = 40.53839999999991
# Row 136
$ws.Cells.Item(136, 8).Value = 2899999.8
$ws.Cells.Item(136, 9).Value = 4832341
$ws.Cells.Item(136, 11).Value = 14497023
$ws.Cells.Item(136, 13).Value = -14494473

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 374.94736
$ws.Cells.Item(5, 9).Value = 363
$ws.Cells.Item(5, 11).Value = 1089
$ws.Cells.Item(5, 13).Value = -977
# Row 88
$ws.Cells.Item(88, 8).Value = 4430.2856
$ws.Cells.Item(88, 9).Value = 2507
$ws.Cells.Item(88, 10).Value = 5199.6
$ws.Cells.Item(88, 11).Value = 7521
$ws.Cells.Item(88, 12).Value = 15598.8
$ws.Cells.Item(88, 13).Value = -7093
$ws.Cells.Item(88, 14).Value = -16454.8
# Row 91
$ws.Cells.Item(91, 8).Value = 4430.2856
$ws.Cells.Item(91, 9).Value = 2507
$ws.Cells.Item(91, 10).Value = 5199.6
$ws.Cells.Item(91, 11).Value = 7521
$ws.Cells.Item(91, 12).Value = 15598.8
$ws.Cells.Item(91, 13).Value = -6039
$ws.Cells.Item(91, 14).Value = -18562.8
# Row 107
$ws.Cells.Item(107, 8).Value = 1860.3
$ws.Cells.Item(107, 10).Value = 2054.1428
$ws.Cells.Item(107, 12).Value = 6162.428400000001
$ws.Cells.Item(107, 14).Value = -10002.4284
# Row 131
$ws.Cells.Item(131, 8).Value = 8487461
$ws.Cells.Item(131, 10).Value = 13551.446
$ws.Cells.Item(131, 12).Value = 40654.338
$ws.Cells.Item(131, 14).Value = -50734.338
# Row 132
$ws.Cells.Item(132, 8).Value = 1666.2858
$ws.Cells.Item(132, 9).Value = 1299.5
$ws.Cells.Item(132, 10).Value = 1813
$ws.Cells.Item(132, 11).Value = 11695.5
$ws.Cells.Item(132, 12).Value = 16317
$ws.Cells.Item(132, 13).Value = -9165.5
$ws.Cells.Item(132, 14).Value = -21377
# Row 135
$ws.Cells.Item(135, 8).Value = 374.94736
$ws.Cells.Item(135, 9).Value = 363
$ws.Cells.Item(135, 11).Value = 3267
$ws.Cells.Item(135, 13).Value = -732

$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Cells.Item(44, 8).Value = 9745
$ws.Cells.Item(44, 9).Value = 9490
$ws.Cells.Item(44, 11).Value = 9490
$ws.Cells.Item(44, 13).Value = -8894
# Row 49
$ws.Cells.Item(49, 8).Value = 24999.5
$ws.Cells.Item(49, 10).Value = 24999.5
$ws.Cells.Item(49, 12).Value = 24999.5
$ws.Cells.Item(49, 14).Value = -25367.5
# Row 55
$ws.Cells.Item(55, 8).Value = 10852
$ws.Cells.Item(55, 9).Value = 13130
$ws.Cells.Item(55, 10).Value = 9333.333000000001
$ws.Cells.Item(55, 11).Value = 13130
$ws.Cells.Item(55, 12).Value = 9333.333000000001
$ws.Cells.Item(55, 13).Value = -12803
$ws.Cells.Item(55, 14).Value = -9987.333000000001
# Row 70
$ws.Cells.Item(70, 8).Value = 4866.6665
$ws.Cells.Item(70, 10).Value = 4333.3335
$ws.Cells.Item(70, 12).Value = 4333.3335
$ws.Cells.Item(70, 14).Value = -4873.3335
# Row 73
$ws.Cells.Item(73, 8).Value = 4866.6665
$ws.Cells.Item(73, 10).Value = 4333.3335
$ws.Cells.Item(73, 12).Value = 4333.3335
$ws.Cells.Item(73, 14).Value = -6205.3335
# Row 110
$ws.Cells.Item(110, 8).Value = 99687
$ws.Cells.Item(110, 10).Value = 99687
$ws.Cells.Item(110, 12).Value = 99687
$ws.Cells.Item(110, 14).Value = -107867

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Cells.Item(82, 8).Value = 3253.2
$ws.Cells.Item(82, 9).Value = 1583.3334
$ws.Cells.Item(82, 11).Value = 1583.3334
$ws.Cells.Item(82, 13).Value = -1222.3334
# Row 85
$ws.Cells.Item(85, 8).Value = 3253.2
$ws.Cells.Item(85, 9).Value = 1583.3334
$ws.Cells.Item(85, 11).Value = 1583.3334
$ws.Cells.Item(85, 13).Value = -335.3334
# Row 93
$ws.Cells.Item(93, 8).Value = 1527.5
$ws.Cells.Item(93, 9).Value = 1075
$ws.Cells.Item(93, 10).Value = 1980
$ws.Cells.Item(93, 11).Value = 1075
$ws.Cells.Item(93, 12).Value = 1980
$ws.Cells.Item(93, 13).Value = 173
$ws.Cells.Item(93, 14).Value = -4476
# Row 136
$ws.Cells.Item(136, 8).Value = 2765.5715
$ws.Cells.Item(136, 9).Value = 1956.3715
$ws.Cells.Item(136, 11).Value = 5869.1145
$ws.Cells.Item(136, 13).Value = -3319.1145

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 90903.28999999999
$ws.Cells.Item(122, 9).Value = 105729.39
$ws.Cells.Item(122, 10).Value = 1946.6666
$ws.Cells.Item(122, 11).Value = 317188.17
$ws.Cells.Item(122, 12).Value = 5839.9998
$ws.Cells.Item(122, 13).Value = -314738.17
$ws.Cells.Item(122, 14).Value = -10739.9998
# Row 133
$ws.Cells.Item(133, 8).Value = 38939.332
$ws.Cells.Item(133, 10).Value = 38939.332
$ws.Cells.Item(133, 12).Value = 38939.332
$ws.Cells.Item(133, 14).Value = -49059.332
